# Apply cryptocurrency price/volume updates from the "cryptos list" refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "0.3940", "11.30") -
# force the Text number format so Excel keeps the exact literal instead of
# normalizing/stripping trailing zeros.
$priceCells = @("D2","D3","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D39","D40","D41","D43","D44","D45","D46","D47","D49","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "22.020.81"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "1.552.75"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "289.94"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("D7").Value = "0.3940"
$ws.Range("E7").Value = "  +3.44%  "
$ws.Range("D8").Value = "0.3218"
$ws.Range("E8").Value = "  -2.99%  "
$ws.Range("D9").Value = "43.87"
$ws.Range("E9").Value = "  -1.83%  "
$ws.Range("D10").Value = "0.07233"
$ws.Range("E10").Value = "  -2.26%  "
$ws.Range("D11").Value = "1.075"
$ws.Range("E11").Value = "  -6.26%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "5.679"
$ws.Range("E13").Value = "  -2.82%  "
$ws.Range("D14").Value = "18.73"
$ws.Range("E14").Value = "  -7.39%  "
$ws.Range("D15").Value = "0.00001126"
$ws.Range("E15").Value = "  +4.78%  "
$ws.Range("D16").Value = "6.623"
$ws.Range("E16").Value = "  -1.88%  "
$ws.Range("D17").Value = "1.555.10"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("D18").Value = "0.06569"
$ws.Range("E18").Value = "  -1.29%  "
$ws.Range("D19").Value = "83.42"
$ws.Range("E19").Value = "  -3.52%  "
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "6.266"
$ws.Range("E21").Value = "  -2.11%  "
$ws.Range("D22").Value = "15.49"
$ws.Range("E22").Value = "  -4.10%  "
$ws.Range("D23").Value = "11.30"
$ws.Range("E23").Value = "  -3.67%  "
$ws.Range("D24").Value = "22.034.57"
$ws.Range("E24").Value = "  -0.77%  "
$ws.Range("D25").Value = "2.371"
$ws.Range("E25").Value = "  +4.36%  "
$ws.Range("D26").Value = "2.409"
$ws.Range("E26").Value = "  -5.87%  "
$ws.Range("D27").Value = "148.79"
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("D28").Value = "18.52"
$ws.Range("E28").Value = "  -3.82%  "
$ws.Range("D29").Value = "4.880"
$ws.Range("D30").Value = "1.730.84"
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").Value = "118.67"
$ws.Range("D32").Value = "0.9725"
$ws.Range("E32").Value = "  -10.84%  "
$ws.Range("D33").Value = "5.838"
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("D34").Value = "0.08329"
$ws.Range("E34").Value = "  +1.41%  "
$ws.Range("D35").Value = "9.099"
$ws.Range("E35").Value = "  -2.49%  "
$ws.Range("D36").Value = "1.599"
$ws.Range("E36").Value = "  -16.16%  "
$ws.Range("D37").Value = "0.02261"
$ws.Range("E37").Value = "  -3.12%  "
$ws.Range("E38").Value = "  -4.06%  "
$ws.Range("D39").Value = "0.05986"
$ws.Range("E39").Value = "  -5.37%  "
$ws.Range("D40").Value = "1.207"
$ws.Range("E40").Value = "  -2.06%  "
$ws.Range("D41").Value = "0.2030"
$ws.Range("E41").Value = "  -6.11%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Value = "10.66"
$ws.Range("E43").Value = "  -3.02%  "
$ws.Range("D44").Value = "0.5800"
$ws.Range("E44").Value = "  -4.40%  "
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").Value = "3.741"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "12.96"
$ws.Range("E46").Value = "  -5.62%  "
$ws.Range("D47").Value = "0.5565"
$ws.Range("E47").Value = "  -5.21%  "
$ws.Range("E48").Value = "  -3.55%  "
$ws.Range("D49").Value = "117.88"
$ws.Range("E49").Value = "  -3.63%  "
$ws.Range("E50").Value = "  -3.84%  "
$ws.Range("D51").Value = "0.06815"
